$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 6
$ws.Range("C3").Value = 8
$ws.Range("C5").Value = 7
$ws.Range("C7").Value = 7
$ws.Range("C9").Value = 7
$ws.Range("C11").Value = 5
$ws.Range("C15").Value = 5
$ws.Range("C18").Value = 11
